$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column O (bud_bank_location): "rootstock" -> "basal_buds" ---
$rootstockRows = @(3,4,6,7,8,9,14,15,16,20)
foreach ($r in $rootstockRows) {
    $ws.Range("O$r").Value = "basal_buds"
    $ws.Range("P$r").ClearContents()
}

# --- Column R (fire_ephemeral): "obligate_fire_ephemeral" -> "fire_ephemeral_obligate" ---
$obligateRows = @(37,39,41)
foreach ($r in $obligateRows) {
    $ws.Range("R$r").Value = "fire_ephemeral_obligate"
}

# --- Column U (life_history_ephemeral) removed entirely ---
$ws.Range("U1").ClearContents()
$ws.Range("U37").ClearContents()
$ws.Range("U39").ClearContents()
$ws.Range("U41").ClearContents()

# --- Move the active selection, matching the author's last cursor position ---
[void]$ws.Range("M2").Select()
